# Weekly update: a new Jengibre price record is inserted at the top of the
# historical series (row 76, right under the existing header/data block),
# pushing the previously existing rows 76-109 down to rows 77-110.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 76; Excel shifts rows 76:109 down to 77:110
# and carries the existing formatting (e.g. the date number format on
# column D) into the freshly inserted row.
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(76, 1).Value = 8
$ws.Cells.Item(76, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76, 3).Value = "Coquimbo"
$ws.Cells.Item(76, 4).Value = 44981
$ws.Cells.Item(76, 5).Value = 4
$ws.Cells.Item(76, 6).Value = 100114007
$ws.Cells.Item(76, 7).Value = "Jengibre"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 400
$ws.Cells.Item(76, 11).Value = 19000
$ws.Cells.Item(76, 12).Value = 20000
$ws.Cells.Item(76, 13).Value = 19500
$ws.Cells.Item(76, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(76, 15).Value = "Perú"
$ws.Cells.Item(76, 16).Value = 1500
$ws.Cells.Item(76, 17).Value = 13
$ws.Cells.Item(76, 18).Value = "Hortaliza"
